# counting.xlsx — "Add files via upload" edit
#
# 1. The Image column (A) held bare filenames like "counting_0.png"; the
#    images were moved into a "counting/" subfolder, so every such value
#    gets a "counting/" prefix.
# 2. One question's wording changed from "cushions" to "pillows".
# 3. One reference answer (row 11, the "typical spoons" question) changed
#    from a single numeric answer to the text "4 and 5 both ok".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $imgCell = $ws.Cells.Item($r, 1)
    $imgVal = $imgCell.Text
    if ($imgVal -match '^counting_\d+\.png$') {
        $imgCell.Value = "counting/$imgVal"
    }

    $qCell = $ws.Cells.Item($r, 2)
    if ($qCell.Text -eq 'How many cushions are depicted in the picture?') {
        $qCell.Value = 'How many pillows are depicted in the picture?'
    }
}

$ws.Range("C11").Value = '4 and 5 both ok'
